$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold the current date, shifting
# all existing header columns one place to the right (REF NUM -> B, etc.)
$ws.Range("A1").EntireColumn.Insert()

$ws.Range("A1").Value = "DATE"
$ws.Range("A1").Style = $ws.Range("B1").Style
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
